$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header updates ---
$ws.Range("A8").Value = "Volume 32   Number  48"
$ws.Range("C9").Value = "Report Covering the Week  11/24/2025  Through  11/30/2025"

# --- Data table updates ---
$ws.Range("D14").Value = 1
$ws.Range("D14").NumberFormat = '#,##0'
$ws.Range("E14").Value = -100
$ws.Range("E14").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("G14").Value = 2
$ws.Range("J14").Value = 9
$ws.Range("K14").Value = -66.666666666666
$ws.Range("L14").Value = -66.666666666666
$ws.Range("D15").Value = 1
$ws.Range("D15").NumberFormat = '#,##0'
$ws.Range("E15").Value = -100
$ws.Range("E15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = -50
$ws.Range("J15").Value = 36
$ws.Range("K15").Value = 19.444444444444
$ws.Range("D16").Value = 13
$ws.Range("E16").Value = -84.615384615384
$ws.Range("F16").Value = 8
$ws.Range("G16").Value = 21
$ws.Range("H16").Value = -61.904761904761
$ws.Range("I16").Value = 242
$ws.Range("J16").Value = 264
$ws.Range("K16").Value = -8.333333333333
$ws.Range("L16").Value = -20.394736842105
$ws.Range("M16").Value = -9.022556390977
$ws.Range("N16").Value = -69.444444444444
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 9
$ws.Range("E17").Value = -55.555555555555
$ws.Range("F17").Value = 25
$ws.Range("G17").Value = 40
$ws.Range("H17").Value = -37.5
$ws.Range("I17").Value = 361
$ws.Range("J17").Value = 420
$ws.Range("K17").Value = -14.047619047619
$ws.Range("L17").Value = -19.777777777777
$ws.Range("M17").Value = 0.277777777777
$ws.Range("N17").Value = -24.476987447698
$ws.Range("C18").Value = 2
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 16
$ws.Range("G18").Value = 7
$ws.Range("H18").Value = 128.571428571429
$ws.Range("I18").Value = 196
$ws.Range("J18").Value = 171
$ws.Range("K18").Value = 14.619883040935
$ws.Range("L18").Value = 18.072289156626
$ws.Range("M18").Value = -2
$ws.Range("N18").Value = -79.835390946502
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 33
$ws.Range("G19").Value = 34
$ws.Range("H19").Value = -2.941176470588
$ws.Range("I19").Value = 347
$ws.Range("J19").Value = 404
$ws.Range("K19").Value = -14.108910891089
$ws.Range("L19").Value = -16.385542168674
$ws.Range("M19").Value = 39.357429718875
$ws.Range("N19").Value = -14.950980392156
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = 0
$ws.Range("G20").Value = 17
$ws.Range("H20").Value = -52.941176470588
$ws.Range("I20").Value = 155
$ws.Range("J20").Value = 166
$ws.Range("K20").Value = -6.626506024096
$ws.Range("L20").Value = -44.444444444444
$ws.Range("M20").Value = 115.277777777778
$ws.Range("N20").Value = -61.728395061728
$ws.Range("C21").Value = 20
$ws.Range("D21").Value = 38
$ws.Range("E21").Value = -47.368421052631
$ws.Range("F21").Value = 91
$ws.Range("G21").Value = 123
$ws.Range("H21").Value = -26.016260162601
$ws.Range("I21").Value = 1347
$ws.Range("J21").Value = 1470
$ws.Range("K21").Value = -8.367346938775
$ws.Range("L21").Value = -18.065693430656
$ws.Range("M21").Value = 15.029888983774
$ws.Range("N21").Value = -56.909788867562
$ws.Range("D22").Value = 2
$ws.Range("G22").Value = 7
$ws.Range("H22").Value = -85.714285714285
$ws.Range("J22").Value = 24
$ws.Range("K22").Value = -37.5
$ws.Range("D23").Value = 2
$ws.Range("D23").NumberFormat = '#,##0'
$ws.Range("E23").Value = -100
$ws.Range("E23").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("G23").Value = 3
$ws.Range("H23").Value = -33.333333333333
$ws.Range("J23").Value = 26
$ws.Range("K23").Value = -57.692307692307
$ws.Range("M23").Value = -8.333333333333
$ws.Range("C24").Value = 17
$ws.Range("D24").Value = 10
$ws.Range("E24").Value = 70
$ws.Range("F24").Value = 43
$ws.Range("H24").Value = -37.681159420289
$ws.Range("I24").Value = 740
$ws.Range("J24").Value = 742
$ws.Range("K24").Value = -0.269541778975
$ws.Range("L24").Value = -11.589008363201
$ws.Range("M24").Value = 20.717781402936
$ws.Range("C25").Value = 5
$ws.Range("C25").NumberFormat = '#,##0'
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 9
$ws.Range("G25").Value = 14
$ws.Range("H25").Value = -35.714285714285
$ws.Range("I25").Value = 160
$ws.Range("J25").Value = 203
$ws.Range("K25").Value = -21.182266009852
$ws.Range("L25").Value = -36.507936507936
$ws.Range("C26").Value = 11
$ws.Range("D26").Value = 13
$ws.Range("E26").Value = -15.384615384615
$ws.Range("F26").Value = 48
$ws.Range("G26").Value = 56
$ws.Range("H26").Value = -14.285714285714
$ws.Range("I26").Value = 597
$ws.Range("J26").Value = 612
$ws.Range("K26").Value = -2.450980392156
$ws.Range("L26").Value = 12.218045112782
$ws.Range("M26").Value = -2.610114192495
$ws.Range("D27").Value = 2
$ws.Range("D27").NumberFormat = '#,##0'
$ws.Range("E27").Value = -100
$ws.Range("E27").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = -75
$ws.Range("J27").Value = 60
$ws.Range("K27").Value = 11.666666666666
$ws.Range("L27").Value = 59.523809523809
$ws.Range("C28").Value = 1
$ws.Range("C28").NumberFormat = '#,##0'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "***.*"
$ws.Range("E28").NumberFormat = "General"
$ws.Range("F28").Value = 1
$ws.Range("G28").Value = 7
$ws.Range("H28").Value = -85.714285714285
$ws.Range("I28").Value = 82
$ws.Range("K28").Value = -30.508474576271
$ws.Range("L28").Value = -10.869565217391
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "***.*"
$ws.Range("E29").NumberFormat = "General"
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = 0
$ws.Range("N29").Value = -88.392857142857
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "***.*"
$ws.Range("E30").NumberFormat = "General"
$ws.Range("G30").Value = 1
$ws.Range("H30").Value = 0
$ws.Range("N30").Value = -88.172043010752
